# Natmi following Dr Hou advice
# The Lgi3-Flot1 ligand-receptor pair table is recomputed: the sending/target
# cluster set grows from {FAPs} to {FAPs, sCs}, so the 4-row table (one row
# per target cluster, for sending cluster "FAPs") becomes an 8-row table
# (4 target clusters x 2 sending clusters: "FAPs" and "sCs").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#          Ligand-expressing cells, Ligand detection rate,
#          Ligand average/total expression value,
#          Ligand derived specificity of average/total expression value,
#          Receptor-expressing cells, Receptor detection rate,
#          Receptor average/total expression value,
#          Receptor derived specificity of average/total expression value,
#          Edge average/total expression weight,
#          Edge average/total expression derived specificity
$data = @(
    @("FAPs", "Lgi3", "Flot1", "ECs", 3, 1, 0.6167776666666667, 1.850333, 0.7836323164322263, 0.7836323164322262, 3, 1, 16.28388233333333, 48.851647, 0.3486636287871168, 0.3486636287871168, 10.04353494982789, 90.39181454845099, 0.2732240870821142, 0.2732240870821142),
    @("FAPs", "Lgi3", "Flot1", "FAPs", 3, 1, 0.6167776666666667, 1.850333, 0.7836323164322263, 0.7836323164322262, 3, 1, 14.832885, 44.498655, 0.3175954851316676, 0.3175954851316677, 9.148592200234999, 82.33732980211499, 0.2488780857021454, 0.2488780857021454),
    @("FAPs", "Lgi3", "Flot1", "M2", 3, 1, 0.6167776666666667, 1.850333, 0.7836323164322263, 0.7836323164322262, 3, 1, 10.02049633333333, 30.061489, 0.214554646263697, 0.214554646263697, 6.180418347315222, 55.623765125837, 0.1681319544529178, 0.1681319544529178),
    @("FAPs", "Lgi3", "Flot1", "sCs", 3, 1, 0.6167776666666667, 1.850333, 0.7836323164322263, 0.7836323164322262, 3, 1, 5.566438666666667, 16.699316, 0.1191862398175185, 0.1191862398175185, 3.433255052469778, 30.899295472228, 0.09339818919504887, 0.09339818919504887),
    @("sCs", "Lgi3", "Flot1", "ECs", 2, 0.6666666666666666, 0.1702976666666667, 0.510893, 0.2163676835677737, 0.2163676835677737, 3, 1, 16.28388233333333, 48.851647, 0.3486636287871168, 0.3486636287871168, 2.773107165641223, 24.957964490771, 0.0754395417050026, 0.0754395417050026),
    @("sCs", "Lgi3", "Flot1", "FAPs", 2, 0.6666666666666666, 0.1702976666666667, 0.510893, 0.2163676835677737, 0.2163676835677737, 3, 1, 14.832885, 44.498655, 0.3175954851316676, 0.3175954851316677, 2.526005705435, 22.734051348915, 0.06871739942952224, 0.06871739942952225),
    @("sCs", "Lgi3", "Flot1", "M2", 2, 0.6666666666666666, 0.1702976666666667, 0.510893, 0.2163676835677737, 0.2163676835677737, 3, 1, 10.02049633333333, 30.061489, 0.214554646263697, 0.214554646263697, 1.706467144408556, 15.358204299677, 0.0464226918107792, 0.04642269181077921),
    @("sCs", "Lgi3", "Flot1", "sCs", 2, 0.6666666666666666, 0.1702976666666667, 0.510893, 0.2163676835677737, 0.2163676835677737, 3, 1, 5.566438666666667, 16.699316, 0.1191862398175185, 0.1191862398175185, 0.9479515165764445, 8.531563649188001, 0.02578805062246963, 0.02578805062246964)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

